# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ 3 = 1391; 5 = 110; 7 = 11808; 8 = 4416; 10 = 45; 13 = 2557; 15 = 157; 17 = 5128; 21 = 11370; 22 = 11330; 24 = 49 }
    "全部类型" = @{ 3 = 1391; 5 = 110; 7 = 11808; 8 = 4416; 10 = 45; 13 = 2557; 16 = 157; 18 = 5128; 22 = 11370; 23 = 11330; 25 = 49 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $ws.Range("F$rowNum").Value = $rows[$rowNum]
    }
}
